# Allocate revenue from EV charger deployment to ISIC codes in shares
# defined in an input data file.
#
# Adds a new "EV Charger Revenue Share by Recipient ISIC Code" (EVCRSbRIC)
# variable row to the "Key to Variables" sheet, alphabetically sorted right
# after the existing "EVCC" (EV Charger Cost) row (row 193), pushing the
# rest of the "trans" section down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new blank row at row 193 (shifts rows 193-214 down to 194-215).
$ws.Rows.Item(193).Insert()

# Populate the new row: Top Level Folder, Acronym, Meaning, Importance.
$ws.Range("A193").Value = "trans"
$ws.Range("B193").Value = "EVCRSbRIC"
$ws.Range("C193").Value = "EV Charger Revenue Share by Recipient ISIC Code"
$ws.Range("F193").Value = "low"
